$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$images = @(
  "flower/flower012.jpg",
  "flower/flower000.jpg",
  "flower/flower028.jpg",
  "flower/flower008.jpg",
  "flower/flower026.jpg",
  "dog/dog029.jpg",
  "flower/flower024.jpg",
  "flower/flower025.jpg",
  "dog/dog012.jpg",
  "flower/flower015.jpg",
  "flower/flower020.jpg",
  "dog/dog009.jpg",
  "flower/flower014.jpg",
  "flower/flower027.jpg",
  "dog/dog005.jpg",
  "dog/dog020.jpg",
  "dog/dog014.jpg",
  "flower/flower013.jpg",
  "dog/dog003.jpg",
  "dog/dog001.jpg",
  "dog/dog028.jpg",
  "dog/dog007.jpg",
  "dog/dog025.jpg",
  "dog/dog030.jpg",
  "flower/flower003.jpg",
  "dog/dog015.jpg",
  "dog/dog023.jpg",
  "dog/dog022.jpg",
  "dog/dog010.jpg",
  "flower/flower018.jpg",
  "flower/flower029.jpg",
  "flower/flower023.jpg"
)

$words = @(
  "ärgern",
  "achten",
  "legen",
  "danken",
  "angeln",
  "heben",
  "hassen",
  "wehen",
  "stoppen",
  "parken",
  "sparen",
  "rechnen",
  "kennen",
  "lernen",
  "quellen",
  "ändern",
  "lügen",
  "leeren",
  "küssen",
  "bergen",
  "deuten",
  "zielen",
  "münzen",
  "streifen",
  "tollen",
  "spüren",
  "spenden",
  "süßen",
  "trotzen",
  "proben",
  "fließen",
  "prüfen"
)

$ws.Range("B2:D33").Clear()

for ($i = 0; $i -lt $images.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $images[$i]
}

for ($i = 0; $i -lt $words.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

for ($i = 0; $i -lt $images.Count; $i++) {
    $row = $i + 2
    $category = $images[$i].Split("/")[0]
    $ws.Cells.Item($row, 4).Value = $category
}
